$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgf2"
$ws.Cells.Item(2, 3).Value = "Fgfr4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.747119
$ws.Cells.Item(2, 8).Value = 2.241357
$ws.Cells.Item(2, 9).Value = 0.03096954854571248
$ws.Cells.Item(2, 10).Value = 0.03096954854571248
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.132884
$ws.Cells.Item(2, 14).Value = 0.398652
$ws.Cells.Item(2, 15).Value = 0.01195569974366677
$ws.Cells.Item(2, 16).Value = 0.01195569974366677
$ws.Cells.Item(2, 17).Value = 0.099280161196
$ws.Cells.Item(2, 18).Value = 0.8935214507639999
$ws.Cells.Item(2, 19).Value = 0.0003702626236094501
$ws.Cells.Item(2, 20).Value = 0.0003702626236094501

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgf2"
$ws.Cells.Item(3, 3).Value = "Fgfr4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.747119
$ws.Cells.Item(3, 8).Value = 2.241357
$ws.Cells.Item(3, 9).Value = 0.03096954854571248
$ws.Cells.Item(3, 10).Value = 0.03096954854571248
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.07352966666666667
$ws.Cells.Item(3, 14).Value = 0.220589
$ws.Cells.Item(3, 15).Value = 0.006615533976389704
$ws.Cells.Item(3, 16).Value = 0.006615533976389703
$ws.Cells.Item(3, 17).Value = 0.05493541103033334
$ws.Cells.Item(3, 18).Value = 0.494418699273
$ws.Cells.Item(3, 19).Value = 0.0002048801006376113
$ws.Cells.Item(3, 20).Value = 0.0002048801006376112

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgf2"
$ws.Cells.Item(4, 3).Value = "Fgfr4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.747119
$ws.Cells.Item(4, 8).Value = 2.241357
$ws.Cells.Item(4, 9).Value = 0.03096954854571248
$ws.Cells.Item(4, 10).Value = 0.03096954854571248
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 10.908285
$ws.Cells.Item(4, 14).Value = 32.724855
$ws.Cells.Item(4, 15).Value = 0.9814287662799436
$ws.Cells.Item(4, 16).Value = 0.9814287662799435
$ws.Cells.Item(4, 17).Value = 8.149786980915
$ws.Cells.Item(4, 18).Value = 73.348082828235
$ws.Cells.Item(4, 19).Value = 0.03039440582146542
$ws.Cells.Item(4, 20).Value = 0.03039440582146542

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fgf2"
$ws.Cells.Item(5, 3).Value = "Fgfr4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 19.74619233333334
$ws.Cells.Item(5, 8).Value = 59.23857700000001
$ws.Cells.Item(5, 9).Value = 0.8185184181638298
$ws.Cells.Item(5, 10).Value = 0.8185184181638298
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.132884
$ws.Cells.Item(5, 14).Value = 0.398652
$ws.Cells.Item(5, 15).Value = 0.01195569974366677
$ws.Cells.Item(5, 16).Value = 0.01195569974366677
$ws.Cells.Item(5, 17).Value = 2.623953022022667
$ws.Cells.Item(5, 18).Value = 23.615577198204
$ws.Cells.Item(5, 19).Value = 0.009785960442227825
$ws.Cells.Item(5, 20).Value = 0.009785960442227825

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fgf2"
$ws.Cells.Item(6, 3).Value = "Fgfr4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 19.74619233333334
$ws.Cells.Item(6, 8).Value = 59.23857700000001
$ws.Cells.Item(6, 9).Value = 0.8185184181638298
$ws.Cells.Item(6, 10).Value = 0.8185184181638298
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.07352966666666667
$ws.Cells.Item(6, 14).Value = 0.220589
$ws.Cells.Item(6, 15).Value = 0.006615533976389704
$ws.Cells.Item(6, 16).Value = 0.006615533976389703
$ws.Cells.Item(6, 17).Value = 1.451930940205889
$ws.Cells.Item(6, 18).Value = 13.067378461853
$ws.Cells.Item(6, 19).Value = 0.005414936405663572
$ws.Cells.Item(6, 20).Value = 0.005414936405663571

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fgf2"
$ws.Cells.Item(7, 3).Value = "Fgfr4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 19.74619233333334
$ws.Cells.Item(7, 8).Value = 59.23857700000001
$ws.Cells.Item(7, 9).Value = 0.8185184181638298
$ws.Cells.Item(7, 10).Value = 0.8185184181638298
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 10.908285
$ws.Cells.Item(7, 14).Value = 32.724855
$ws.Cells.Item(7, 15).Value = 0.9814287662799436
$ws.Cells.Item(7, 16).Value = 0.9814287662799435
$ws.Cells.Item(7, 17).Value = 215.397093636815
$ws.Cells.Item(7, 18).Value = 1938.573842731335
$ws.Cells.Item(7, 19).Value = 0.8033175213159384
$ws.Cells.Item(7, 20).Value = 0.8033175213159383

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Fgf2"
$ws.Cells.Item(8, 3).Value = "Fgfr4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 3.630999
$ws.Cells.Item(8, 8).Value = 10.892997
$ws.Cells.Item(8, 9).Value = 0.1505120332904577
$ws.Cells.Item(8, 10).Value = 0.1505120332904577
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.132884
$ws.Cells.Item(8, 14).Value = 0.398652
$ws.Cells.Item(8, 15).Value = 0.01195569974366677
$ws.Cells.Item(8, 16).Value = 0.01195569974366677
$ws.Cells.Item(8, 17).Value = 0.482501671116
$ws.Cells.Item(8, 18).Value = 4.342515040044
$ws.Cells.Item(8, 19).Value = 0.001799476677829488
$ws.Cells.Item(8, 20).Value = 0.001799476677829489

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Fgf2"
$ws.Cells.Item(9, 3).Value = "Fgfr4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 3.630999
$ws.Cells.Item(9, 8).Value = 10.892997
$ws.Cells.Item(9, 9).Value = 0.1505120332904577
$ws.Cells.Item(9, 10).Value = 0.1505120332904577
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.07352966666666667
$ws.Cells.Item(9, 14).Value = 0.220589
$ws.Cells.Item(9, 15).Value = 0.006615533976389704
$ws.Cells.Item(9, 16).Value = 0.006615533976389703
$ws.Cells.Item(9, 17).Value = 0.266986146137
$ws.Cells.Item(9, 18).Value = 2.402875315233
$ws.Cells.Item(9, 19).Value = 0.000995717470088521
$ws.Cells.Item(9, 20).Value = 0.000995717470088521

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Fgf2"
$ws.Cells.Item(10, 3).Value = "Fgfr4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.630999
$ws.Cells.Item(10, 8).Value = 10.892997
$ws.Cells.Item(10, 9).Value = 0.1505120332904577
$ws.Cells.Item(10, 10).Value = 0.1505120332904577
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 10.908285
$ws.Cells.Item(10, 14).Value = 32.724855
$ws.Cells.Item(10, 15).Value = 0.9814287662799436
$ws.Cells.Item(10, 16).Value = 0.9814287662799435
$ws.Cells.Item(10, 17).Value = 39.607971926715
$ws.Cells.Item(10, 18).Value = 356.4717473404349
$ws.Cells.Item(10, 19).Value = 0.1477168391425397
$ws.Cells.Item(10, 20).Value = 0.1477168391425397
